# edit.ps1 -- applies the "Added GitHub Link to my file" commit to the
# Group Project Module 9 document via Word COM interop.
#
# Summary of changes (derived from the canonical OOXML diff):
#   1. "Module " + "9 "                              -> merge into "Module 9 "
#   2. "Ricardo Orlando" + ", Mo" + "nica Jones..."   -> merge first two runs
#   3. "Due: 0" + "5" + "/0" + "4" + "/2023"          -> merge into "Due: 05/04/2023"
#   4. "Employee has one " + "position"               -> merge into one run
#   5. "Bacchus has monthly " + "inventory"           -> merge into one run
#   6. Append a new "GitHub Link: <url>" bold line at the end of the document.

$d = $word.ActiveDocument

$wdReplaceAll = 2

# 1. "Module " / "9 " -> "Module 9 " (single merged run)
$d.Content.Find.Execute("Module 9 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Module 9 ", $wdReplaceAll) | Out-Null

# 2. "Ricardo Orlando" + ", Mo" -> "Ricardo Orlando, Mo" (merge into the run
#    that is followed by "nica Jones, Donnell Perkins")
$d.Content.Find.Execute("Ricardo Orlando, Mo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ricardo Orlando, Mo", $wdReplaceAll) | Out-Null

# 3. "Due: 0" + "5" + "/0" + "4" + "/2023" -> "Due: 05/04/2023"
$d.Content.Find.Execute("Due: 05/04/2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Due: 05/04/2023", $wdReplaceAll) | Out-Null

# 4. "Employee has one " + "position" -> "Employee has one position"
$d.Content.Find.Execute("Employee has one position", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Employee has one position", $wdReplaceAll) | Out-Null

# 5. "Bacchus has monthly " + "inventory" -> "Bacchus has monthly inventory"
$d.Content.Find.Execute("Bacchus has monthly inventory", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Bacchus has monthly inventory", $wdReplaceAll) | Out-Null

# 6. Append the GitHub Link line into the final (empty, bold 28pt) paragraph.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)   # wdCollapseEnd

$r1 = $r.Duplicate
$r1.Text = "GitHub Link: "
$r1.Font.Bold = $true
$r1.Font.Size = 14

$r2 = $d.Paragraphs.Last.Range
$r2.Collapse(0)
$r2.Text = "https://github.com/donnellperkins/csd-310.git"
$r2.Font.Bold = $true
$r2.Font.Size = 14
